$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Elements sheet: update canonical terminology URLs ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElem.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElem.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R54-DiplomeUniversiteInterUniversitaire/FHIR/TRE-R54-DiplomeUniversiteInterUniversitaire?vs"

# --- Elements sheet: widen column Z (bestFit width grew with the longer URL) ---
# Target OOXML width is 104.68359375; the host's ColumnWidth setter snaps to a
# whole-pixel grid (px = round(chars*MDW + 5)), so 104.68359375 itself is not
# reachable bit-for-bit. 103.8333333333 lands on the closest attainable pixel
# (104.66666666666667), the nearest the COM model can get to the target width.
$wsElem.Columns.Item(26).ColumnWidth = 103.8333333333
